# Apply "Generate Report for Handback" edits to localization-status.xlsx

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# "Ready for handoff" -> "Handed back: in sync with en-US" for both zh-cn and de-de status columns
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen status columns E & F to fit new, longer text
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet ---
# Refresh handback datetime
$zhcn.Range("K2").Value = "2016-08-12 06:44:20"
$zhcn.Range("K3").Value = "2016-08-12 06:44:20"

# Clear stale "version not latest" error details (now in sync, no error)
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# Widen status column C, narrow now-unused Error Detail column P
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet ---
# Refresh handback datetime
$dede.Range("K2").Value = "2016-08-12 06:44:29"
$dede.Range("K3").Value = "2016-08-12 06:44:29"

# Clear stale "version not latest" error details (now in sync, no error)
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

# Widen status column C, narrow now-unused Error Detail column P
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
